$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells (coin names / URLs / volume percentages) - safe to set directly.
$ws.Range('E2').Value = '  +1.01%  '
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('E6').Value = '  +5.32%  '
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('E10').Value = '  +5.25%  '
$ws.Range('E11').Value = '  +3.76%  '
$ws.Range('E12').Value = '  +2.49%  '
$ws.Range('E13').Value = '  +3.56%  '
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('E15').Value = '  +1.54%  '
$ws.Range('E16').Value = '  +2.27%  '
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('E22').Value = '  +2.64%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('B24').Value = 'Avalanche'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('E24').Value = '  +0.30%  '
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('E27').Value = '  +2.25%  '
$ws.Range('E28').Value = '  +1.23%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('E32').Value = '  +1.58%  '
$ws.Range('E33').Value = '  -4.59%  '
$ws.Range('E34').Value = '  +4.24%  '
$ws.Range('E35').Value = '  +4.95%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('E37').Value = '  +7.13%  '
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('E39').Value = '  +0.86%  '
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  +3.76%  '
$ws.Range('E43').Value = '  +3.55%  '
$ws.Range('E44').Value = '  +7.12%  '
$ws.Range('E45').Value = '  +2.39%  '
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('E48').Value = '  +0.55%  '
$ws.Range('E49').Value = '  +3.09%  '
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('E51').Value = '  -0.08%  '

# Numeric-looking price strings must stay text: force Text format, assign, then restore Normal style
# so no stray "s" style attribute is left behind on the cell.
$priceCells = @{
    'D2' = '26.982.06'
    'D3' = '1.676.95'
    'D5' = '216.01'
    'D7' = '0.999'
    'D10' = '20.40'
    'D12' = '1.912.63'
    'D13' = '1.698.97'
    'D17' = '27.012.46'
    'D18' = '233.22'
    'D19' = '7.84'
    'D23' = '2.22'
    'D24' = '9.23'
    'D25' = '145.57'
    'D26' = '7.15'
    'D28' = '16.00'
    'D29' = '1.00'
    'D33' = '1.458.01'
    'D34' = '3.17'
    'D38' = '0.570'
    'D42' = '2.31'
    'D43' = '65.96'
    'D44' = '0.972'
    'D45' = '1.820.14'
    'D46' = '0.780'
    'D47' = '90.62'
    'D51' = '7.60'
}
foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
    $cell.Style = "Normal"
}

Write-Host "cryptos list updated"
